# Insert 2 new data rows at row 275 (pushing existing rows 275-372 down to 277-374)
# and populate them with the new reported prices.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("275:276").Insert()

# New row 275
$ws.Range("A275").Value = 10
$ws.Range("B275").Value = "Vega Modelo de Temuco"
$ws.Range("C275").Value = "La Araucanía"
$ws.Range("D275").Value = 44468
$ws.Range("E275").Value = 9
$ws.Range("F275").Value = 100112003
$ws.Range("G275").Value = "Ajo"
$ws.Range("H275").Value = "Chino"
$ws.Range("I275").Value = "Primera"
$ws.Range("J275").Value = 300
$ws.Range("K275").Value = 17000
$ws.Range("L275").Value = 18000
$ws.Range("M275").Value = 17667
$ws.Range("N275").Value = "`$/caja 10 kilos"
$ws.Range("O275").Value = "China"
$ws.Range("P275").Value = 1767
$ws.Range("Q275").Value = 10
$ws.Range("R275").Value = "Hortaliza"

# New row 276
$ws.Range("A276").Value = 10
$ws.Range("B276").Value = "Vega Modelo de Temuco"
$ws.Range("C276").Value = "La Araucanía"
$ws.Range("D276").Value = 44468
$ws.Range("E276").Value = 9
$ws.Range("F276").Value = 100112003
$ws.Range("G276").Value = "Ajo"
$ws.Range("H276").Value = "Chino"
$ws.Range("I276").Value = "Primera"
$ws.Range("J276").Value = 100
$ws.Range("K276").Value = 19000
$ws.Range("L276").Value = 19000
$ws.Range("M276").Value = 19000
$ws.Range("N276").Value = "`$/malla 10 kilos"
$ws.Range("O276").Value = "China"
$ws.Range("P276").Value = 1900
$ws.Range("Q276").Value = 10
$ws.Range("R276").Value = "Hortaliza"
